$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(19, 8).Value = 378
$ws.Cells.Item(19, 9).Value = 432.3846
$ws.Cells.Item(19, 10).Value = 333.8125
$ws.Cells.Item(19, 11).Value = 432.3846
$ws.Cells.Item(19, 12).Value = 333.8125
$ws.Cells.Item(19, 13).Value = -257.3846
$ws.Cells.Item(19, 14).Value = -683.8125
$ws.Cells.Item(132, 8).Value = 5752820.5
$ws.Cells.Item(132, 9).Value = 6413811
$ws.Cells.Item(132, 10).Value = 24234.334
$ws.Cells.Item(132, 11).Value = 19241433
$ws.Cells.Item(132, 12).Value = 72703.00199999999
$ws.Cells.Item(132, 13).Value = -19238903
$ws.Cells.Item(132, 14).Value = -77763.00199999999
$ws.Cells.Item(136, 8).Value = 40173.332
$ws.Cells.Item(136, 10).Value = 40173.332
$ws.Cells.Item(136, 12).Value = 40173.332
$ws.Cells.Item(136, 14).Value = -50373.332

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 8488.674000000001
$ws.Cells.Item(32, 9).Value = 6070.845
$ws.Cells.Item(32, 11).Value = 6070.845
$ws.Cells.Item(32, 13).Value = -5783.845
$ws.Cells.Item(61, 8).Value = 1364.4231
$ws.Cells.Item(61, 9).Value = 1247.2174
$ws.Cells.Item(61, 10).Value = 2263
$ws.Cells.Item(61, 11).Value = 1247.2174
$ws.Cells.Item(61, 12).Value = 2263
$ws.Cells.Item(61, 13).Value = -1035.2174
$ws.Cells.Item(61, 14).Value = -2687
$ws.Cells.Item(110, 8).Value = 428
$ws.Cells.Item(110, 9).Value = 428
$ws.Cells.Item(110, 11).Value = 428
$ws.Cells.Item(110, 13).Value = 1617
$ws.Cells.Item(132, 8).Value = 2251.2144
$ws.Cells.Item(132, 9).Value = 1475.9756
$ws.Cells.Item(132, 10).Value = 4370.2
$ws.Cells.Item(132, 11).Value = 4427.9268
$ws.Cells.Item(132, 12).Value = 13110.6
$ws.Cells.Item(132, 13).Value = -1897.9268
$ws.Cells.Item(132, 14).Value = -18170.6
$ws.Cells.Item(136, 8).Value = 1364.4231
$ws.Cells.Item(136, 9).Value = 1247.2174
$ws.Cells.Item(136, 10).Value = 2263
$ws.Cells.Item(136, 11).Value = 3741.6522
$ws.Cells.Item(136, 12).Value = 6789
$ws.Cells.Item(136, 13).Value = -1191.6522
$ws.Cells.Item(136, 14).Value = -11889

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 10870236
$ws.Cells.Item(94, 9).Value = 19231388
$ws.Cells.Item(94, 10).Value = 737
$ws.Cells.Item(94, 11).Value = 19231388
$ws.Cells.Item(94, 12).Value = 737
$ws.Cells.Item(94, 13).Value = -19230937
$ws.Cells.Item(94, 14).Value = -1639
$ws.Cells.Item(105, 8).Value = 500050000
$ws.Cells.Item(105, 9).Value = 500050000
$ws.Cells.Item(105, 11).Value = 500050000
$ws.Cells.Item(105, 13).Value = -500048253

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 45455596
$ws.Cells.Item(16, 9).Value = 62501068
$ws.Cells.Item(16, 11).Value = 62501068
$ws.Cells.Item(16, 13).Value = -62500781
$ws.Cells.Item(58, 8).Value = 3964.0227
$ws.Cells.Item(58, 9).Value = 1076.1111
$ws.Cells.Item(58, 10).Value = 8550.706
$ws.Cells.Item(58, 11).Value = 1076.1111
$ws.Cells.Item(58, 12).Value = 8550.706
$ws.Cells.Item(58, 13).Value = -873.1111000000001
$ws.Cells.Item(58, 14).Value = -8956.706
$ws.Cells.Item(94, 8).Value = 1922.9412
$ws.Cells.Item(94, 9).Value = 1453.1428
$ws.Cells.Item(94, 10).Value = 2251.8
$ws.Cells.Item(94, 11).Value = 1453.1428
$ws.Cells.Item(94, 12).Value = 2251.8
$ws.Cells.Item(94, 13).Value = -1002.1428
$ws.Cells.Item(94, 14).Value = -3153.8
$ws.Cells.Item(107, 8).Value = 606.7
$ws.Cells.Item(107, 9).Value = 277.13333
$ws.Cells.Item(107, 11).Value = 277.13333
$ws.Cells.Item(107, 13).Value = 1642.86667
$ws.Cells.Item(113, 8).Value = 45455596
$ws.Cells.Item(113, 9).Value = 62501068
$ws.Cells.Item(113, 11).Value = 62501068
$ws.Cells.Item(113, 13).Value = -62498898
$ws.Cells.Item(122, 8).Value = 1213.6666
$ws.Cells.Item(122, 9).Value = 1104
$ws.Cells.Item(122, 10).Value = 1433
$ws.Cells.Item(122, 11).Value = 3312
$ws.Cells.Item(122, 12).Value = 4299
$ws.Cells.Item(122, 13).Value = -862
$ws.Cells.Item(122, 14).Value = -9199
$ws.Cells.Item(134, 8).Value = 1400.4286
$ws.Cells.Item(134, 9).Value = 1382.2142
$ws.Cells.Item(134, 11).Value = 4146.642599999999
$ws.Cells.Item(134, 13).Value = -1611.642599999999
$ws.Cells.Item(136, 8).Value = 3964.0227
$ws.Cells.Item(136, 9).Value = 1076.1111
$ws.Cells.Item(136, 10).Value = 8550.706
$ws.Cells.Item(136, 11).Value = 3228.3333
$ws.Cells.Item(136, 12).Value = 25652.118
$ws.Cells.Item(136, 13).Value = -678.3333000000002
$ws.Cells.Item(136, 14).Value = -30752.118
$ws.Cells.Item(141, 8).Value = 480933.06
$ws.Cells.Item(141, 10).Value = 519344.16
$ws.Cells.Item(141, 12).Value = 519344.16
$ws.Cells.Item(141, 14).Value = -529704.1599999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 4568557
$ws.Cells.Item(4, 10).Value = 4772533.5
$ws.Cells.Item(4, 12).Value = 14317600.5
$ws.Cells.Item(4, 14).Value = -14317824.5
$ws.Cells.Item(40, 8).Value = 231.25
$ws.Cells.Item(40, 10).Value = 344.44446
$ws.Cells.Item(40, 12).Value = 1377.77784
$ws.Cells.Item(40, 14).Value = -1515.77784
$ws.Cells.Item(131, 8).Value = 27068438
$ws.Cells.Item(131, 9).Value = 62500510
$ws.Cells.Item(131, 11).Value = 187501530
$ws.Cells.Item(131, 13).Value = -187496490
$ws.Cells.Item(132, 8).Value = 1191.4286
$ws.Cells.Item(132, 9).Value = 959.375
$ws.Cells.Item(132, 10).Value = 1500.8334
$ws.Cells.Item(132, 11).Value = 8634.375
$ws.Cells.Item(132, 12).Value = 13507.5006
$ws.Cells.Item(132, 13).Value = -6104.375
$ws.Cells.Item(132, 14).Value = -18567.5006

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 1808.775
$ws.Cells.Item(102, 9).Value = 1363.6154
$ws.Cells.Item(102, 11).Value = 1363.6154
$ws.Cells.Item(102, 13).Value = 258.3846000000001
$ws.Cells.Item(122, 8).Value = 1060
$ws.Cells.Item(122, 9).Value = 1013.3333
$ws.Cells.Item(122, 10).Value = 1200
$ws.Cells.Item(122, 11).Value = 3039.9999
$ws.Cells.Item(122, 12).Value = 3600
$ws.Cells.Item(122, 13).Value = -589.9998999999998
$ws.Cells.Item(122, 14).Value = -8500
$ws.Cells.Item(126, 8).Value = 2089.1333
$ws.Cells.Item(126, 9).Value = 1703.8
$ws.Cells.Item(126, 10).Value = 2859.8
$ws.Cells.Item(126, 11).Value = 5111.4
$ws.Cells.Item(126, 12).Value = 8579.400000000001
$ws.Cells.Item(126, 13).Value = -2641.4
$ws.Cells.Item(126, 14).Value = -13519.4
$ws.Cells.Item(132, 8).Value = 2483.7083
$ws.Cells.Item(132, 9).Value = 2316.111
$ws.Cells.Item(132, 10).Value = 2986.5
$ws.Cells.Item(132, 11).Value = 6948.333
$ws.Cells.Item(132, 12).Value = 8959.5
$ws.Cells.Item(132, 13).Value = -4418.333
$ws.Cells.Item(132, 14).Value = -14019.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 2960.5557
$ws.Cells.Item(7, 9).Value = 2877.1428
$ws.Cells.Item(7, 10).Value = 3252.5
$ws.Cells.Item(7, 11).Value = 2877.1428
$ws.Cells.Item(7, 12).Value = 3252.5
$ws.Cells.Item(7, 13).Value = -2765.1428
$ws.Cells.Item(7, 14).Value = -3476.5
$ws.Cells.Item(40, 8).Value = 4820.7856
$ws.Cells.Item(40, 9).Value = 2280.0908
$ws.Cells.Item(40, 10).Value = 14136.667
$ws.Cells.Item(40, 11).Value = 2280.0908
$ws.Cells.Item(40, 12).Value = 2280.0908
$ws.Cells.Item(40, 13).Value = -2144.0908
$ws.Cells.Item(40, 14).Value = -14408.667
$ws.Cells.Item(93, 8).Value = 950
$ws.Cells.Item(93, 9).Value = 933.3333
$ws.Cells.Item(93, 10).Value = 1000
$ws.Cells.Item(93, 11).Value = 933.3333
$ws.Cells.Item(93, 12).Value = 1000
$ws.Cells.Item(93, 13).Value = 314.6667
$ws.Cells.Item(93, 14).Value = -3496
$ws.Cells.Item(126, 8).Value = 2960.5557
$ws.Cells.Item(126, 9).Value = 2877.1428
$ws.Cells.Item(126, 10).Value = 3252.5
$ws.Cells.Item(126, 11).Value = 8631.428400000001
$ws.Cells.Item(126, 12).Value = 9757.5
$ws.Cells.Item(126, 13).Value = -6161.428400000001
$ws.Cells.Item(126, 14).Value = -14697.5
$ws.Cells.Item(133, 8).Value = 52299.668
$ws.Cells.Item(133, 10).Value = 52299.668
$ws.Cells.Item(133, 12).Value = 52299.668
$ws.Cells.Item(133, 14).Value = -57359.668
$ws.Cells.Item(136, 8).Value = 1670.4706
$ws.Cells.Item(136, 9).Value = 1276.3846
$ws.Cells.Item(136, 11).Value = 3829.1538
$ws.Cells.Item(136, 13).Value = -1279.1538

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(46, 8).Value = 38197.4
$ws.Cells.Item(46, 10).Value = 38197.4
$ws.Cells.Item(46, 12).Value = 38197.4
$ws.Cells.Item(46, 14).Value = -38659.4
$ws.Cells.Item(113, 8).Value = 406
$ws.Cells.Item(113, 9).Value = 245
$ws.Cells.Item(113, 10).Value = 567
$ws.Cells.Item(113, 11).Value = 735
$ws.Cells.Item(113, 12).Value = 1701
$ws.Cells.Item(113, 13).Value = 1435
$ws.Cells.Item(113, 14).Value = -6041
$ws.Cells.Item(122, 8).Value = 10002213
$ws.Cells.Item(122, 10).Value = 2221
$ws.Cells.Item(122, 12).Value = 6663
$ws.Cells.Item(122, 14).Value = -11563
$ws.Cells.Item(126, 8).Value = 55556380
$ws.Cells.Item(126, 9).Value = 62500770
$ws.Cells.Item(126, 10).Value = 1252.5
$ws.Cells.Item(126, 11).Value = 187502310
$ws.Cells.Item(126, 12).Value = 3757.5
$ws.Cells.Item(126, 13).Value = -187499840
$ws.Cells.Item(126, 14).Value = -8697.5
$ws.Cells.Item(134, 8).Value = 38197.4
$ws.Cells.Item(134, 10).Value = 38197.4
$ws.Cells.Item(134, 12).Value = 114592.2
$ws.Cells.Item(134, 14).Value = -119662.2
